$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q2 and R2 to nearest integer
$ws.Range("Q2").Value = [math]::Round($ws.Range("Q2").Value(), 0)
$ws.Range("R2").Value = [math]::Round($ws.Range("R2").Value(), 0)

# Clear Z2 and AB2 entirely (they had "00:00" inline strings)
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
